# "Some print statement change" — rework the youtubedata header row:
#   A: date -> Title
#   B: Title -> url
#   C: FileType (unchanged)
#   D: url column removed entirely
# and drop the sample data row (row 2) that was only there to exercise
# the old "date" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-point the header labels (C1 "FileType" is already correct and stays put).
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "url"

# The old D column ("url") is no longer needed now that column B holds it.
$ws.Columns("D:D").Delete()

# Drop the sample data row beneath the header.
$ws.Rows("2:2").Delete()
